# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets to the latest scraped values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 805.1111
$ws.Range("I33").Value = 805.1111
$ws.Range("K33").Value = 805.1111
$ws.Range("M33").Value = -576.1111
$ws.Range("H80").Value = 6227.8
$ws.Range("I80").Value = 814.5714
$ws.Range("J80").Value = 9142.615
$ws.Range("K80").Value = 2443.7142
$ws.Range("L80").Value = 27427.845
$ws.Range("M80").Value = -1445.7142
$ws.Range("N80").Value = -29423.845
$ws.Range("H83").Value = 6227.8
$ws.Range("I83").Value = 814.5714
$ws.Range("J83").Value = 9142.615
$ws.Range("K83").Value = 7331.1426
$ws.Range("L83").Value = 82283.535
$ws.Range("M83").Value = -2339.1426
$ws.Range("N83").Value = -92267.535
$ws.Range("H116").Value = 9526347
$ws.Range("I116").Value = 200000000
$ws.Range("J116").Value = 2664
$ws.Range("K116").Value = 200000000
$ws.Range("L116").Value = 2664
$ws.Range("M116").Value = -199996558
$ws.Range("N116").Value = -9548
$ws.Range("H127").Value = 1076.1052
$ws.Range("I127").Value = 538.55554
$ws.Range("J127").Value = 1559.9
$ws.Range("K127").Value = 1615.66662
$ws.Range("L127").Value = 4679.700000000001
$ws.Range("M127").Value = 3344.33338
$ws.Range("N127").Value = -14599.7
$ws.Range("H129").Value = 1125.3778
$ws.Range("I129").Value = 451.2
$ws.Range("J129").Value = 1318
$ws.Range("K129").Value = 1353.6
$ws.Range("L129").Value = 3954
$ws.Range("M129").Value = 3646.4
$ws.Range("N129").Value = -13954
$ws.Range("H135").Value = 824.8333
$ws.Range("I135").Value = 777.3333
$ws.Range("J135").Value = 1157.3334
$ws.Range("K135").Value = 6995.9997
$ws.Range("L135").Value = 10416.0006
$ws.Range("M135").Value = -4460.9997
$ws.Range("N135").Value = -15486.0006
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21472.104
$ws.Range("I32").Value = 22409.72
$ws.Range("J32").Value = 15612
$ws.Range("K32").Value = 22409.72
$ws.Range("L32").Value = 15612
$ws.Range("M32").Value = -22122.72
$ws.Range("N32").Value = -16186
$ws.Range("H45").Value = 993.4091
$ws.Range("I45").Value = 971.17645
$ws.Range("J45").Value = 1069
$ws.Range("K45").Value = 971.17645
$ws.Range("L45").Value = 1069
$ws.Range("M45").Value = -594.17645
$ws.Range("N45").Value = -1823
$ws.Range("I46").Value = 2138
$ws.Range("J46").Value = 5500
$ws.Range("K46").Value = 2138
$ws.Range("L46").Value = 5500
$ws.Range("M46").Value = -1819
$ws.Range("N46").Value = -6138
$ws.Range("H88").Value = 2913.2856
$ws.Range("I88").Value = 2447.6667
$ws.Range("K88").Value = 2447.6667
$ws.Range("M88").Value = -2041.6667
$ws.Range("H91").Value = 2913.2856
$ws.Range("I91").Value = 2447.6667
$ws.Range("K91").Value = 2447.6667
$ws.Range("M91").Value = -1043.6667
$ws.Range("H114").Value = 40049
$ws.Range("J114").Value = 40049
$ws.Range("L114").Value = 40049
$ws.Range("N114").Value = -48727
$ws.Range("H119").Value = 27954.2
$ws.Range("J119").Value = 27954.2
$ws.Range("L119").Value = 27954.2
$ws.Range("N119").Value = -37630.2
$ws.Range("H122").Value = 3341.2273
$ws.Range("I122").Value = 3478.6667
$ws.Range("J122").Value = 3046.7144
$ws.Range("K122").Value = 10436.0001
$ws.Range("L122").Value = 9140.143199999999
$ws.Range("M122").Value = -7986.000100000001
$ws.Range("N122").Value = -14040.1432
$ws.Range("H123").Value = 24157.857
$ws.Range("J123").Value = 24157.857
$ws.Range("L123").Value = 24157.857
$ws.Range("N123").Value = -33957.857
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1422.2307
$ws.Range("I99").Value = 849
$ws.Range("J99").Value = 3333
$ws.Range("K99").Value = 849
$ws.Range("L99").Value = 3333
$ws.Range("M99").Value = 649
$ws.Range("N99").Value = -6329
$ws.Range("H134").Value = 2399.535
$ws.Range("I134").Value = 2086.5312
$ws.Range("K134").Value = 6259.5936
$ws.Range("M134").Value = -3724.5936
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2080.5625
$ws.Range("I31").Value = 1172.9565
$ws.Range("K31").Value = 1172.9565
$ws.Range("M31").Value = -877.9565
$ws.Range("H34").Value = 2080.5625
$ws.Range("I34").Value = 1172.9565
$ws.Range("K34").Value = 1172.9565
$ws.Range("M34").Value = -970.9565
$ws.Range("H107").Value = 373.34784
$ws.Range("I107").Value = 310.8
$ws.Range("J107").Value = 490.625
$ws.Range("K107").Value = 310.8
$ws.Range("L107").Value = 490.625
$ws.Range("M107").Value = 1609.2
$ws.Range("N107").Value = -4330.625
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 821.25
$ws.Range("I122").Value = 522.5
$ws.Range("J122").Value = 1120
$ws.Range("K122").Value = 4702.5
$ws.Range("L122").Value = 10080
$ws.Range("M122").Value = -2252.5
$ws.Range("N122").Value = -14980
$ws.Range("H132").Value = 1621.2759
$ws.Range("I132").Value = 1176.4615
$ws.Range("J132").Value = 1982.6875
$ws.Range("K132").Value = 10588.1535
$ws.Range("L132").Value = 17844.1875
$ws.Range("M132").Value = -8058.153499999999
$ws.Range("N132").Value = -22904.1875
$ws.Range("H133").Value = 3643.7827
$ws.Range("J133").Value = 6582
$ws.Range("L133").Value = 19746
$ws.Range("N133").Value = -29866
$ws.Range("H134").Value = 3455.0527
$ws.Range("I134").Value = 1800.9131
$ws.Range("J134").Value = 5991.4
$ws.Range("K134").Value = 5402.7393
$ws.Range("L134").Value = 17974.2
$ws.Range("M134").Value = -332.7393000000002
$ws.Range("N134").Value = -28114.2
$ws.Range("H136").Value = 4287.609
$ws.Range("I136").Value = 1018.5
$ws.Range("J136").Value = 6031.1333
$ws.Range("K136").Value = 3055.5
$ws.Range("L136").Value = 18093.3999
$ws.Range("M136").Value = 2044.5
$ws.Range("N136").Value = -28293.3999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 44125.332
$ws.Range("I97").Value = 50483.08
$ws.Range("J97").Value = 2800
$ws.Range("K97").Value = 50483.08
$ws.Range("L97").Value = 2800
$ws.Range("M97").Value = -49987.08
$ws.Range("N97").Value = -3792
$ws.Range("H122").Value = 3672.1667
$ws.Range("I122").Value = 2975
$ws.Range("J122").Value = 4020.75
$ws.Range("K122").Value = 8925
$ws.Range("L122").Value = 12062.25
$ws.Range("M122").Value = -6475
$ws.Range("N122").Value = -16962.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1266.2727
$ws.Range("I16").Value = 1055.7
$ws.Range("J16").Value = 3372
$ws.Range("K16").Value = 1055.7
$ws.Range("L16").Value = 3372
$ws.Range("M16").Value = -885.7
$ws.Range("N16").Value = -3712
$ws.Range("H22").Value = 902
$ws.Range("I22").Value = 1057.3334
$ws.Range("J22").Value = 817.2727
$ws.Range("K22").Value = 1057.3334
$ws.Range("L22").Value = 817.2727
$ws.Range("M22").Value = -762.3334
$ws.Range("N22").Value = -1407.2727
$ws.Range("H27").Value = 902
$ws.Range("I27").Value = 1057.3334
$ws.Range("J27").Value = 817.2727
$ws.Range("K27").Value = 1057.3334
$ws.Range("L27").Value = 817.2727
$ws.Range("M27").Value = -950.3334
$ws.Range("N27").Value = -1031.2727
$ws.Range("H61").Value = 30172
$ws.Range("I61").Value = 34367.332
$ws.Range("K61").Value = 34367.332
$ws.Range("M61").Value = -34165.332
$ws.Range("H113").Value = 30172
$ws.Range("I113").Value = 34367.332
$ws.Range("K113").Value = 34367.332
$ws.Range("M113").Value = -32197.332
$ws.Range("H122").Value = 25006000
$ws.Range("I122").Value = 22732546
$ws.Range("K122").Value = 68197638
$ws.Range("M122").Value = -68195188
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H62").Value = 4900
$ws.Range("J62").Value = 4900
$ws.Range("L62").Value = 4900
$ws.Range("N62").Value = -6148
$ws.Range("H65").Value = 4900
$ws.Range("J65").Value = 4900
$ws.Range("L65").Value = 24500
$ws.Range("N65").Value = -30740
$ws.Range("H119").Value = 160474.25
$ws.Range("J119").Value = 160474.25
$ws.Range("L119").Value = 160474.25
$ws.Range("N119").Value = -170150.25
$ws.Range("H122").Value = 46298812
$ws.Range("J122").Value = 3219.3333
$ws.Range("L122").Value = 9657.999899999999
$ws.Range("N122").Value = -14557.9999
